$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B6").Value = 8
$ws.Range("B7").Value = 8
$ws.Range("B8").Value = 14
$ws.Range("B9").Value = 11

$ws.Range("B9").Select()
